$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Fix typo in a company name (stray space before 司)
$ws.Range("B7").Value = "潤泰創新國際股份有限公司"

# Insert a new "property_category" column before the existing "date" column (H),
# shifting date / legislator_name / legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"

$lastRow = 23
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
